{"js": "// Fix mistake in one word: \"Use the video in our \" -> \"Use our \"\n// (also relocates the auto-tracked \"_GoBack\" bookmark to the edit point,\n// matching what Word itself does when text is edited/retyped).\n\n// 1) Drop the old \"_GoBack\" bookmark whose stale location was inside the\n//    earlier empty paragraph.\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\n// 2) Find the text we need to remove (\"the video in \") inside the\n//    \"Use the video in our \" run.\nconst results = context.document.body.search(\"the video in \", { matchCase: true, matchWholeWord: false });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error('Could not find \"the video in \" in the document body.');\n}\n\nconst toRemove = results.items[0];\n\n// 3) Re-create \"_GoBack\" as a collapsed bookmark right at the point of the\n//    edit (i.e. right before the text we are about to delete) \u2014 this is\n//    what splits the run into \"Use \" / \"our \" around the bookmark.\nconst editPoint = toRemove.getRange(\"Start\");\neditPoint.insertBookmark(\"_GoBack\");\nawait context.sync();\n\n// 4) Finally remove the now-unwanted text itself.\ntoRemove.insertText(\"\", Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# Fix mistake in one word: \"Use the video in our \" -> \"Use our \"\n# (also relocates the auto-tracked \"_GoBack\" bookmark to the edit point,\n# matching what Word itself does when text is edited/retyped).\n\n$d = $word.ActiveDocument\n\n# 1) Drop the old \"_GoBack\" bookmark whose stale location was inside the\n#    earlier empty paragraph.\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks.Item(\"_GoBack\").Delete()\n}\n\n# 2) Find the text we need to remove (\"the video in \") inside the\n#    \"Use the video in our \" sentence.\n$find = $d.Content.Find\n$find.Text = \"the video in \"\n$found = $find.Execute()\nif (-not $found) {\n    throw 'Could not find \"the video in \" in the document.'\n}\n$toRemove = $find.Parent\n\n# 3) Re-create \"_GoBack\" as a collapsed bookmark right at the point of the\n#    edit (i.e. right before the text we are about to delete) -- this is\n#    what splits the run into \"Use \" / \"our \" around the bookmark.\n$editPoint = $d.Range($toRemove.Start, $toRemove.Start)\n$d.Bookmarks.Add(\"_GoBack\", $editPoint)\n\n# 4) Finally remove the now-unwanted text itself.\n$toRemove.Delete()\n"}
